# Automation the third test in Class name OpenAccountTest.Java
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet for the OpenAccountTest scenario, after the existing sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "asmaa dawood"
$ws2.Range("B2").Value = "Dollar"

$ws2.Range("B2").Select()

# Make the new sheet the active tab
$ws2.Activate()
